$wb = $excel.ActiveWorkbook

# Add a new worksheet named "before" at the end
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "before"

# Populate the new sheet with the two new loop-listener example strings in column B
$newSheet.Range("B1").Value = '<jt:for var="x" start="1" end="5" onLoopProcessed="${blockShadingLoopListener}" copyRight="true">${x}</jt:for>'
$newSheet.Range("B2").Value = '<jt:for var="x" start="1" end="1" onLoopProcessed="${blockShadingLoopListener}">The above count, using ${x}, should have 3 replaced!</jt:for>'

# Move the new sheet to after the last existing sheet (implInstance) so it becomes the last tab
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)
